$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New category data entered on rows 42-44 ---
# C42 = "Half", E42 = comment about delete screen
$ws.Range("C42").Value = "Half"
$ws.Range("E42").Value = "Delete работи, но не се отваря delete screen"

# C43 and C44 = "Yes"
$ws.Range("C43").Value = "Yes"
$ws.Range("C44").Value = "Yes"

# --- Widen column E to fit the new, longer comment ---
$ws.Columns.Item(5).ColumnWidth = 52.25

# --- Update the view: scroll down to the new rows and select C44 ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C44").Select() | Out-Null
